# Update TPM-derived specificity/weight values in the LR-pairs sheet
# (Cxcl13-Cxcr3) to reflect results recomputed with new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("I2").Value = 0.8781048434890718
$ws.Range("J2").Value = 0.8781048434890719
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("S2").Value = 0.01414779397916906
$ws.Range("T2").Value = 0.01414779397916906

# Row 3
$ws.Range("I3").Value = 0.8781048434890718
$ws.Range("J3").Value = 0.8781048434890719
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 9.533166330895556
$ws.Range("R3").Value = 85.79849697806
$ws.Range("S3").Value = 0.5999939335250724
$ws.Range("T3").Value = 0.5999939335250725

# Row 4
$ws.Range("I4").Value = 0.8781048434890718
$ws.Range("J4").Value = 0.8781048434890719
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 4.194049555002223
$ws.Range("R4").Value = 37.74644599502
$ws.Range("S4").Value = 0.2639631159848304
$ws.Range("T4").Value = 0.2639631159848303

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.489861
$ws.Range("H5").Value = 1.469583
$ws.Range("I5").Value = 0.1218951565109281
$ws.Range("J5").Value = 0.1218951565109281
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("Q5").Value = 0.031204635561
$ws.Range("R5").Value = 0.280841720049
$ws.Range("S5").Value = 0.001963942659196415
$ws.Range("T5").Value = 0.001963942659196415

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.489861
$ws.Range("H6").Value = 1.469583
$ws.Range("I6").Value = 0.1218951565109281
$ws.Range("J6").Value = 0.1218951565109281
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 1.323357695343
$ws.Range("R6").Value = 11.910219258087
$ws.Range("S6").Value = 0.08328886348245759
$ws.Range("T6").Value = 0.08328886348245759

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.489861
$ws.Range("H7").Value = 1.469583
$ws.Range("I7").Value = 0.1218951565109281
$ws.Range("J7").Value = 0.1218951565109281
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 9.533166330895556
$ws.Range("R7").Value = 5.239817291079
$ws.Range("S7").Value = 0.03664235036927409
$ws.Range("T7").Value = 0.03664235036927409
